# Refresh the "cryptos" price table (columns D = Price, E = Volume(1h))
# with the latest scraped figures.
#
# Cell values are assigned with a leading apostrophe so Excel treats them
# as literal text instead of auto-converting numeric-looking strings to
# numbers (which would silently drop things like trailing zeros or
# multi-dot formatted prices). ClearFormats() afterwards removes the
# "text quote prefix" cell style that the apostrophe trick leaves behind,
# so the cells keep using the sheet default style exactly as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.613.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -1.51%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.903.79"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.41%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'526.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -2.41%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'143.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -4.89%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.549"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -3.52%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.909.08"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -2.54%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -5.04%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -2.41%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D13").Value = "'3.417.58"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -2.20%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'  +2.56%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'60.591.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -1.66%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'22.57"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -4.64%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.909.12"
$ws.Range("D17").ClearFormats()
$ws.Range("E18").Value = "'  -3.85%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -3.29%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  -3.71%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'350.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -7.88%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -2.18%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'5.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +0.72%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'64.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.76%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  -3.98%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.179"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -4.98%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").Value = "'7.82"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -4.81%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.0₃0852"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -9.56%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  -2.77%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'19.59"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -4.33%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'151.33"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -4.80%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'4.32"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -6.98%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'5.57"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -5.66%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.998"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -6.85%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'1.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -5.85%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'37.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +0.42%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -4.91%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'3.72"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -5.05%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'2.291.95"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -5.03%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.649"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -3.33%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  -1.99%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'20.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -7.46%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'4.95"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -2.82%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'  -3.07%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'10.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -1.02%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0918"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'18.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -7.51%  "
$ws.Range("E51").ClearFormats()
